$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N ("Late"), shifting the existing
# N/O/P ("Late", "heading"/Original, "Outstanding") columns right by one.
$ws.Columns("N").Insert()

# The newly inserted column should keep the same width as column M (its
# left neighbour), matching Excel's default "copy formatting from the left"
# behaviour when inserting a column.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and select cell R9 on it,
# matching the workbook's saved view state.
$ws.Activate() | Out-Null
$ws.Range("R9").Select() | Out-Null
